# Update the COVID-19 "paises" dashboard: refresh case counts for several
# countries and update the "last updated" timestamp. A few countries
# (Costa Rica, Haiti, Islas Malvinas) overtake their neighbours in the
# ranking (the sheet is kept sorted by total cases, column B, descending),
# which is why their rows change identity/order as well as value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 22:31"

# --- Straight data refreshes (no change in rank/row) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 7435608
$ws.Cells.Item(4, 3).Value = 29462
$ws.Cells.Item(4, 4).Value = 4679688
$ws.Cells.Item(4, 5).Value = 2544455
$ws.Cells.Item(4, 7).Value = 680
$ws.Cells.Item(4, 8).Value = 211465

# Row 5: India
$ws.Cells.Item(5, 2).Value = 6310267
$ws.Cells.Item(5, 3).Value = 86748
$ws.Cells.Item(5, 4).Value = 5270007
$ws.Cells.Item(5, 5).Value = 941552
$ws.Cells.Item(5, 7).Value = 1179
$ws.Cells.Item(5, 8).Value = 98708

# Row 25: Alemania
$ws.Cells.Item(25, 2).Value = 292896
$ws.Cells.Item(25, 3).Value = 2430
$ws.Cells.Item(25, 5).Value = 27325
$ws.Cells.Item(25, 7).Value = 15
$ws.Cells.Item(25, 8).Value = 9571

# Row 27: Israel
$ws.Cells.Item(27, 2).Value = 245494
$ws.Cells.Item(27, 3).Value = 8568
$ws.Cells.Item(27, 4).Value = 174506
$ws.Cells.Item(27, 5).Value = 69419
$ws.Cells.Item(27, 7).Value = 41
$ws.Cells.Item(27, 8).Value = 1569

# --- Costa Rica overtakes Portugal and Etiopia (rows 52-54) ---
# Row 52 becomes Costa Rica with refreshed numbers
$ws.Cells.Item(52, 1).Value = "Costa Rica"
$ws.Cells.Item(52, 2).Value = 75760
$ws.Cells.Item(52, 3).Value = 1156
$ws.Cells.Item(52, 4).Value = 37841
$ws.Cells.Item(52, 5).Value = 37015
$ws.Cells.Item(52, 7).Value = 24
$ws.Cells.Item(52, 8).Value = 904

# Row 53 becomes Portugal (its own unchanged numbers, shifted down a row)
$ws.Cells.Item(53, 1).Value = "Portugal"
$ws.Cells.Item(53, 2).Value = 75542
$ws.Cells.Item(53, 3).Value = 825
$ws.Cells.Item(53, 4).Value = 48530
$ws.Cells.Item(53, 5).Value = 25041
$ws.Cells.Item(53, 7).Value = 8
$ws.Cells.Item(53, 8).Value = 1971

# Row 54 becomes Etiopia (its own unchanged numbers, shifted down a row)
$ws.Cells.Item(54, 1).Value = "Etiopia"
$ws.Cells.Item(54, 2).Value = 75368
$ws.Cells.Item(54, 3).Value = 784
$ws.Cells.Item(54, 4).Value = 31204
$ws.Cells.Item(54, 5).Value = 42966
$ws.Cells.Item(54, 7).Value = 7
$ws.Cells.Item(54, 8).Value = 1198

# --- Haiti overtakes Gabon (rows 108-109) ---
# Row 108 becomes Haiti with refreshed numbers
$ws.Cells.Item(108, 1).Value = "Haiti"
$ws.Cells.Item(108, 2).Value = 8766
$ws.Cells.Item(108, 3).Value = 26
$ws.Cells.Item(108, 4).Value = 6829
$ws.Cells.Item(108, 5).Value = 1708
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 229

# Row 109 becomes Gabon (its own unchanged numbers, shifted down a row)
$ws.Cells.Item(109, 1).Value = "Gabon"
$ws.Cells.Item(109, 2).Value = 8752
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 7955
$ws.Cells.Item(109, 5).Value = 743
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 54

# --- Islas Malvinas overtakes Montserrat (rows 215-216) ---
# Row 215 becomes Islas Malvinas (its own unchanged numbers)
$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 13
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

# Row 216 becomes Montserrat (its own unchanged numbers)
$ws.Cells.Item(216, 1).Value = "Montserrat"
$ws.Cells.Item(216, 2).Value = 13
$ws.Cells.Item(216, 3).Value = 0
$ws.Cells.Item(216, 4).Value = 12
$ws.Cells.Item(216, 5).Value = 0
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 1
